# Sync attendance_reports: fix "Recorded By" (column G) ordering on the
# "Session Analysis Results" sheet so entries read "System, ..." / the
# gmail.com / admin.com address first, matching the source-of-truth order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Session Analysis Results")

# Map of row -> expected current value -> corrected value for column G
# ("Recorded By"). OldValue is kept alongside NewValue purely for
# documentation/traceability of what each row is being changed from.
$updates = @(
    @{ Row = 2; OldValue = "backup@backdoor.com, System, system"; NewValue = "System, backup@backdoor.com, system" }
    @{ Row = 4; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 5; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 8; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 11; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 17; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 29; OldValue = "backup@backdoor.com, System, system"; NewValue = "System, backup@backdoor.com, system" }
    @{ Row = 31; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 32; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 35; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 38; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 44; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 56; OldValue = "backup@backdoor.com, System, system"; NewValue = "System, backup@backdoor.com, system" }
    @{ Row = 58; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 59; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 62; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 65; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 71; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 83; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 84; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 85; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 90; OldValue = "admin@admin.com, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, admin@admin.com" }
    @{ Row = 96; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 97; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 99; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 109; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 110; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 111; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 116; OldValue = "admin@admin.com, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, admin@admin.com" }
    @{ Row = 122; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 123; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 125; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 135; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 136; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 137; OldValue = "backup@backdoor.com, System"; NewValue = "System, backup@backdoor.com" }
    @{ Row = 142; OldValue = "admin@admin.com, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, admin@admin.com" }
    @{ Row = 148; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 149; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
    @{ Row = 151; OldValue = "System, dnasr281@gmail.com"; NewValue = "dnasr281@gmail.com, System" }
)

foreach ($update in $updates) {
    $ws.Cells.Item($update.Row, 7).Value = $update.NewValue
}

